$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'20.622.88"
$ws.Range('E2').Value = '  +2.46%  '

$ws.Range('D3').Value = "'1.469.91"
$ws.Range('E3').Value = '  +2.69%  '

$ws.Range('D4').Value = "'1.003"
$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').Value = "'0.9806"
$ws.Range('E5').Value = '  -2.06%  '

$ws.Range('D6').Value = "'281.01"
$ws.Range('E6').Value = '  +1.94%  '

$ws.Range('D7').Value = "'0.3734"
$ws.Range('E7').Value = '  +0.79%  '

$ws.Range('D8').Value = "'0.3226"
$ws.Range('E8').Value = '  +4.39%  '

$ws.Range('D9').Value = "'41.87"
$ws.Range('E9').Value = '  +4.46%  '

$ws.Range('D10').Value = "'1.081"
$ws.Range('E10').Value = '  +7.08%  '

$ws.Range('D11').Value = "'0.06805"
$ws.Range('E11').Value = '  +3.31%  '

$ws.Range('D12').Value = "'0.9937"
$ws.Range('E12').Value = '  -0.72%  '

$ws.Range('D13').Value = "'5.674"
$ws.Range('E13').Value = '  +4.35%  '

$ws.Range('D14').Value = "'18.80"
$ws.Range('E14').Value = '  +8.71%  '

$ws.Range('D15').Value = "'6.370"
$ws.Range('E15').Value = '  +3.20%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = "'1.470.72"
$ws.Range('E16').Value = '  +2.74%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = "'0.00001050"
$ws.Range('E17').Value = '  +3.89%  '

$ws.Range('D18').Value = "'0.05827"
$ws.Range('E18').Value = '  -0.29%  '

$ws.Range('B19').Value = 'Litecoin'
$ws.Range('C19').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D19').Value = "'73.61"
$ws.Range('E19').Value = '  -2.33%  '

$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = "'0.9793"
$ws.Range('E20').Value = '  -2.17%  '

$ws.Range('D21').Value = "'5.765"
$ws.Range('E21').Value = '  +0.98%  '

$ws.Range('D22').Value = "'15.16"
$ws.Range('E22').Value = '  +4.29%  '

$ws.Range('D23').Value = "'11.35"
$ws.Range('E23').Value = '  +2.42%  '

$ws.Range('D24').Value = "'2.315"
$ws.Range('E24').Value = '  -0.39%  '

$ws.Range('D25').Value = "'20.671.07"
$ws.Range('E25').Value = '  +2.66%  '

$ws.Range('D26').Value = "'2.377"
$ws.Range('E26').Value = '  +3.16%  '

$ws.Range('D27').Value = "'137.85"
$ws.Range('E27').Value = '  -0.86%  '

$ws.Range('D28').Value = "'17.85"
$ws.Range('E28').Value = '  +4.98%  '

$ws.Range('D29').Value = "'1.629.81"
$ws.Range('E29').Value = '  +2.43%  '

$ws.Range('D30').Value = "'114.84"
$ws.Range('E30').Value = '  +4.89%  '

$ws.Range('D31').Value = "'3.998"
$ws.Range('E31').Value = '  +2.94%  '

$ws.Range('D32').Value = "'5.473"
$ws.Range('E32').Value = '  -0.14%  '

$ws.Range('D33').Value = "'0.8668"
$ws.Range('E33').Value = '  -7.15%  '

$ws.Range('D34').Value = "'0.07924"
$ws.Range('E34').Value = '  +1.97%  '

$ws.Range('D35').Value = "'1.529"
$ws.Range('E35').Value = '  +16.24%  '

$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = "'0.06019"
$ws.Range('E36').Value = '  +6.32%  '

$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = "'5.014"
$ws.Range('E37').Value = '  +4.97%  '

$ws.Range('D38').Value = "'10.98"
$ws.Range('E38').Value = '  -4.86%  '

$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = "'7.985"
$ws.Range('E39').Value = '  -5.21%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.02106"
$ws.Range('E40').Value = '  +3.99%  '

$ws.Range('B41').Value = 'Frax'
$ws.Range('C41').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D41').Value = "'0.9844"
$ws.Range('E41').Value = '  -1.64%  '

$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = "'0.1939"
$ws.Range('E42').Value = '  +1.16%  '

$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = "'1.135"
$ws.Range('E43').Value = '  +1.39%  '

$ws.Range('D44').Value = "'0.5491"
$ws.Range('E44').Value = '  +3.01%  '

$ws.Range('D45').Value = "'12.67"
$ws.Range('E45').Value = '  +3.43%  '

$ws.Range('D46').Value = "'3.596"
$ws.Range('E46').Value = '  +1.01%  '

$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = "'0.5473"
$ws.Range('E47').Value = '  +6.42%  '

$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = "'121.79"
$ws.Range('E48').Value = '  +9.45%  '

$ws.Range('D49').Value = "'1.852"
$ws.Range('E49').Value = '  +3.34%  '

$ws.Range('D50').Value = "'1.071"
$ws.Range('E50').Value = '  +1.74%  '

$ws.Range('D51').Value = "'0.06460"
$ws.Range('E51').Value = '  +3.89%  '
